$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("tasas")

$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.04 = 32934.08 pesos`n✅ 32934.08 pesos = 8.01 = 945.55 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws2.Range("N10").Value = 124.4
$ws2.Range("O10").Value = 4097
$ws2.Range("N12").Value = 4110
$ws2.Range("O12").Value = 118
